$wb = $excel.ActiveWorkbook

# Color constants (BGR ints for Excel COM Interior.Color / matches fills in styles.xml)
$blue   = 13998939   # 5B9BD5
$orange = 49407       # FFC000
$ltblue = 16247773    # DDEBF7
$ltyellow = 13431551  # FFF2CC

# =========================================================================
# Sheet 1 ("MaddenCo Data"): rename tab and re-color the header row so it
# alternates blue / orange fill (Calibri Light, centered, no bold, no top
# alignment) instead of the plain bold unfilled header used previously.
# =========================================================================
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "MaddenCo Data 09-02-2022"

$hdr1Cols = @("A1","B1","C1","D1","E1","F1","G1","H1","I1")
$i = 0
foreach ($addr in $hdr1Cols) {
    $r = $ws1.Range($addr)
    $r.Font.Name = "Calibri Light"
    $r.Font.Size = 12
    $r.Font.Bold = $false
    $r.HorizontalAlignment = -4108   # xlCenter
    $r.VerticalAlignment = -4107     # xlBottom (clears inherited "top")
    $r.Borders.LineStyle = 1
    if ($i % 2 -eq 0) {
        $r.Interior.Color = $blue
    } else {
        $r.Interior.Color = $orange
    }
    $i = $i + 1
}

# =========================================================================
# Sheet 2 ("MaddenCo Counts" -> "MaddenCo Count"): rebuild as a small
# pivot-style table. Row 1 becomes the Employee # header band (no A1
# label), and the previous row 1..4 label/count rows shift down to
# rows 2..5 with their row-label column (A) using the bold, unfilled,
# top-aligned label style.
# =========================================================================
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "MaddenCo Count 09-02-2022"
$ws2.Cells.Clear()

# --- Row 1: Employee # header band, alternating blue / orange fill ---
$empHeaders = @(142,154,161,124,138,149,169)
$headerCols = @("B","C","D","E","F","G","H")
for ($j = 0; $j -lt $headerCols.Length; $j++) {
    $addr = "$($headerCols[$j])1"
    $r = $ws2.Range($addr)
    $r.Value = $empHeaders[$j]
    $r.Font.Name = "Calibri Light"
    $r.Font.Size = 12
    $r.Font.Bold = $false
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4107
    $r.Borders.LineStyle = 1
    if ($j % 2 -eq 0) {
        $r.Interior.Color = $blue
    } else {
        $r.Interior.Color = $orange
    }
}

# --- Rows 2-5: row labels in column A + count grid in B:H ---
$labels = @("Total","Immed","Later","Emails")
$data = @(
    @(5,7,5,14,12,1,3),
    @(2,3,3,4,8,0,0),
    @(3,4,2,10,4,1,3),
    @(3,4,3,3,2,1,3)
)

for ($k = 0; $k -lt $labels.Length; $k++) {
    $rowNum = $k + 2

    $lbl = $ws2.Range("A$rowNum")
    $lbl.Value = $labels[$k]
    $lbl.Font.Name = "Calibri"
    $lbl.Font.Size = 11
    $lbl.Font.Bold = $true
    $lbl.HorizontalAlignment = -4108   # xlCenter
    $lbl.VerticalAlignment = -4160     # xlTop
    $lbl.Borders.LineStyle = 1

    for ($j = 0; $j -lt $headerCols.Length; $j++) {
        $addr = "$($headerCols[$j])$rowNum"
        $r = $ws2.Range($addr)
        $r.Value = $data[$k][$j]
        $r.Font.Name = "Calibri Light"
        $r.Font.Size = 12
        $r.Font.Bold = $false
        $r.HorizontalAlignment = -4108
        $r.VerticalAlignment = -4107
        $r.Borders.LineStyle = 1
        if ($j % 2 -eq 0) {
            $r.Interior.Color = $ltyellow
        } else {
            $r.Interior.Color = $ltblue
        }
    }
}

# Column H has no alternate fill partner beyond G, so force it to match
# the light-blue (odd) band like the original generator's off-by-one loop.
for ($k = 0; $k -lt $labels.Length; $k++) {
    $rowNum = $k + 2
    $ws2.Range("H$rowNum").Interior.Color = $ltblue
}
